$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price/Volume data range as Text so that values such as
# "1.000" or "0.9999" are stored as literal strings (matching the
# original inlineStr cells) instead of being auto-converted to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.652.55'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.873.10'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '248.24'
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '0.4729'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '0.2913'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = '0.06484'
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = '22.07'
$ws.Range("E10").Value = '  +4.83%  '
$ws.Range("D11").Value = '0.07692'
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("D12").Value = '96.61'
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '0.7379'
$ws.Range("E13").Value = '  +2.63%  '
$ws.Range("D14").Value = '1.869.49'
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").Value = '5.166'
$ws.Range("E15").Value = '  +0.49%  '
$ws.Range("D16").Value = '273.33'
$ws.Range("E16").Value = '  -1.19%  '
$ws.Range("D17").Value = '30.655.35'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").Value = '13.35'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("D20").Value = '0.000007528'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '2.117.29'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '5.276'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").Value = '6.189'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").Value = '9.230'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").Value = '164.16'
$ws.Range("E26").Value = '  -0.80%  '
$ws.Range("D27").Value = '18.78'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("D28").Value = '1.911'
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("D29").Value = '0.1000'
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("D30").Value = '1.345'
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").Value = '1.513'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").Value = '4.289'
$ws.Range("E32").Value = '  +0.06%  '
$ws.Range("D33").Value = '4.108'
$ws.Range("E33").Value = '  +1.89%  '
$ws.Range("D34").Value = '0.04807'
$ws.Range("E34").Value = '  +0.61%  '
$ws.Range("D35").Value = '1.121'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = '0.6966'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").Value = '2.717'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '0.01856'
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("D39").Value = '2.752'
$ws.Range("E39").Value = '  +0.33%  '
$ws.Range("D40").Value = '6.252'
$ws.Range("E40").Value = '  -2.37%  '
$ws.Range("D41").Value = '73.26'
$ws.Range("E41").Value = '  +4.26%  '
$ws.Range("D42").Value = '1.976'
$ws.Range("E42").Value = '  +2.99%  '
$ws.Range("D43").Value = '0.4184'
$ws.Range("E43").Value = '  +1.72%  '
$ws.Range("D44").Value = '1.0000'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '0.8343'
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("D46").Value = '101.80'
$ws.Range("E46").Value = '  -0.25%  '
$ws.Range("D47").Value = '9.318'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '35.44'
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").Value = '6.979'
$ws.Range("E49").Value = '  -1.77%  '
$ws.Range("D50").Value = '918.20'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").Value = '0.05652'
$ws.Range("E51").Value = '  +1.30%  '

# Reset the cell style back to Normal now that the text values are
# committed, so no stray formatting is left attached to the cells.
$dataRange.NumberFormat = "General"
$dataRange.Style = "Normal"
